$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 8466.666999999999
$ws.Range("I64").Value = 2999
$ws.Range("J64").Value = 11200.5
$ws.Range("K64").Value = 2999
$ws.Range("L64").Value = 11200.5
$ws.Range("M64").Value = -2751
$ws.Range("N64").Value = -11696.5
$ws.Range("H67").Value = 8466.666999999999
$ws.Range("I67").Value = 2999
$ws.Range("J67").Value = 11200.5
$ws.Range("K67").Value = 2999
$ws.Range("L67").Value = 11200.5
$ws.Range("M67").Value = -2141
$ws.Range("N67").Value = -12916.5
$ws.Range("H101").Value = 1283.6666
$ws.Range("I101").Value = 1003.2857
$ws.Range("J101").Value = 1529
$ws.Range("K101").Value = 3009.8571
$ws.Range("L101").Value = 4587
$ws.Range("M101").Value = -1387.8571
$ws.Range("N101").Value = -7831
$ws.Range("H116").Value = 6490.5454
$ws.Range("I116").Value = 4923
$ws.Range("J116").Value = 7386.2856
$ws.Range("K116").Value = 4923
$ws.Range("L116").Value = 7386.2856
$ws.Range("M116").Value = -1481
$ws.Range("N116").Value = -14270.2856
$ws.Range("H137").Value = 13701537
$ws.Range("I137").Value = 66669210
$ws.Range("J137").Value = 3001.276
$ws.Range("K137").Value = 200007630
$ws.Range("L137").Value = 9003.828
$ws.Range("N137").Value = -14103.828
$ws.Range("H138").Value = 3685.9062
$ws.Range("I138").Value = 2115.3572
$ws.Range("J138").Value = 4125.66
$ws.Range("K138").Value = 6346.071599999999
$ws.Range("L138").Value = 12376.98
$ws.Range("M138").Value = -1206.071599999999
$ws.Range("N138").Value = -22656.98

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5760.7817
$ws.Range("I32").Value = 3445.796
$ws.Range("J32").Value = 24666.5
$ws.Range("K32").Value = 3445.796
$ws.Range("L32").Value = 24666.5
$ws.Range("M32").Value = -3158.796
$ws.Range("N32").Value = -25240.5
$ws.Range("H61").Value = 6634.2354
$ws.Range("I61").Value = 6147.3335
$ws.Range("J61").Value = 7802.8
$ws.Range("K61").Value = 6147.3335
$ws.Range("L61").Value = 7802.8
$ws.Range("M61").Value = -5935.3335
$ws.Range("N61").Value = -8226.799999999999
$ws.Range("H74").Value = 2821.4688
$ws.Range("I74").Value = 1855.12
$ws.Range("J74").Value = 6272.7144
$ws.Range("K74").Value = 1855.12
$ws.Range("L74").Value = 6272.7144
$ws.Range("M74").Value = -981.1199999999999
$ws.Range("N74").Value = -8020.7144
$ws.Range("H77").Value = 2821.4688
$ws.Range("I77").Value = 1855.12
$ws.Range("J77").Value = 6272.7144
$ws.Range("K77").Value = 9275.599999999999
$ws.Range("L77").Value = 31363.572
$ws.Range("M77").Value = -4907.599999999999
$ws.Range("N77").Value = -40099.572
$ws.Range("H97").Value = 984.1667
$ws.Range("I97").Value = 984.1667
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 984.1667
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -488.1667
$ws.Range("N97").ClearContents()
$ws.Range("H102").Value = 3966.6667
$ws.Range("I102").Value = 3966.6667
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 3966.6667
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -2344.6667
$ws.Range("H122").Value = 3739.3635
$ws.Range("I122").Value = 3987.4443
$ws.Range("J122").Value = 3567.6155
$ws.Range("K122").Value = 11962.3329
$ws.Range("L122").Value = 10702.8465
$ws.Range("M122").Value = -9512.332900000001
$ws.Range("N122").Value = -15602.8465
$ws.Range("H132").Value = 3038.1191
$ws.Range("I132").Value = 2070.5806
$ws.Range("J132").Value = 5764.8184
$ws.Range("K132").Value = 6211.7418
$ws.Range("L132").Value = 17294.4552
$ws.Range("M132").Value = -3681.7418
$ws.Range("N132").Value = -22354.4552
$ws.Range("H133").Value = 89990
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 89990
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 89990
$ws.Range("N133").Value = -95050
$ws.Range("H136").Value = 6634.2354
$ws.Range("I136").Value = 6147.3335
$ws.Range("J136").Value = 7802.8
$ws.Range("K136").Value = 18442.0005
$ws.Range("L136").Value = 23408.4
$ws.Range("M136").Value = -15892.0005
$ws.Range("N136").Value = -28508.4
$ws.Range("H137").Value = 69999
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 69999
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 69999
$ws.Range("N137").Value = -80199

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 10000
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 10000
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 10000
$ws.Range("N99").Value = -12996
$ws.Range("M99").ClearContents()
$ws.Range("H105").Value = 12638.6
$ws.Range("I105").Value = 12314.895
$ws.Range("J105").Value = 13197.728
$ws.Range("K105").Value = 12314.895
$ws.Range("L105").Value = 13197.728
$ws.Range("M105").Value = -10567.895
$ws.Range("N105").Value = -16691.728
$ws.Range("H132").Value = 69726
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 69726
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 69726
$ws.Range("N132").Value = -79846
$ws.Range("H134").Value = 2122.72
$ws.Range("I134").Value = 1304.8918
$ws.Range("J134").Value = 4450.385
$ws.Range("K134").Value = 3914.6754
$ws.Range("L134").Value = 13351.155
$ws.Range("M134").Value = -1379.6754
$ws.Range("N134").Value = -18421.155
$ws.Range("H137").Value = 51106.7
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 51106.7
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 51106.7
$ws.Range("N137").Value = -61306.7
$ws.Range("H138").Value = 65436.57
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 65436.57
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 65436.57
$ws.Range("N138").Value = -75716.57000000001
$ws.Range("H140").Value = 66749.60000000001
$ws.Range("I140").Value = 60000
$ws.Range("J140").Value = 67499.55499999999
$ws.Range("K140").Value = 60000
$ws.Range("L140").Value = 67499.55499999999
$ws.Range("M140").Value = -54820
$ws.Range("N140").Value = -77859.55499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2988.8
$ws.Range("I58").Value = 1475.36
$ws.Range("J58").Value = 6772.4
$ws.Range("K58").Value = 1475.36
$ws.Range("L58").Value = 6772.4
$ws.Range("M58").Value = -1272.36
$ws.Range("N58").Value = -7178.4
$ws.Range("H105").Value = 2974.3333
$ws.Range("I105").Value = 1433.5714
$ws.Range("J105").Value = 4322.5
$ws.Range("K105").Value = 1433.5714
$ws.Range("L105").Value = 4322.5
$ws.Range("M105").Value = 313.4286
$ws.Range("N105").Value = -7816.5
$ws.Range("H132").Value = 3174.8914
$ws.Range("I132").Value = 2730.2058
$ws.Range("J132").Value = 4434.8335
$ws.Range("K132").Value = 8190.617400000001
$ws.Range("L132").Value = 13304.5005
$ws.Range("M132").Value = -5660.617400000001
$ws.Range("N132").Value = -18364.5005
$ws.Range("H134").Value = 3332.2727
$ws.Range("I134").Value = 2092.3333
$ws.Range("J134").Value = 5989.2856
$ws.Range("K134").Value = 6276.999899999999
$ws.Range("L134").Value = 17967.8568
$ws.Range("M134").Value = -3741.999899999999
$ws.Range("N134").Value = -23037.8568
$ws.Range("H136").Value = 2988.8
$ws.Range("I136").Value = 1475.36
$ws.Range("J136").Value = 6772.4
$ws.Range("K136").Value = 4426.08
$ws.Range("L136").Value = 20317.2
$ws.Range("M136").Value = -1876.08
$ws.Range("N136").Value = -25417.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3476.2666
$ws.Range("I5").Value = 1168.125
$ws.Range("J5").Value = 6114.143
$ws.Range("K5").Value = 3504.375
$ws.Range("L5").Value = 18342.429
$ws.Range("M5").Value = -3392.375
$ws.Range("N5").Value = -18566.429
$ws.Range("H37").Value = 362727.9
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 362727.9
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 1088183.7
$ws.Range("N37").Value = -1088407.7
$ws.Range("H131").Value = 10306106
$ws.Range("I131").Value = 8929633
$ws.Range("J131").Value = 11439672
$ws.Range("K131").Value = 26788899
$ws.Range("L131").Value = 34319016
$ws.Range("M131").Value = -26783859
$ws.Range("N131").Value = -34329096
$ws.Range("H132").Value = 3601.3215
$ws.Range("I132").Value = 2975.25
$ws.Range("J132").Value = 4436.0835
$ws.Range("K132").Value = 26777.25
$ws.Range("L132").Value = 39924.7515
$ws.Range("M132").Value = -24247.25
$ws.Range("N132").Value = -44984.7515
$ws.Range("H133").Value = 13338128
$ws.Range("I133").Value = 1056.3334
$ws.Range("J133").Value = 20840232
$ws.Range("K133").Value = 3169.0002
$ws.Range("L133").Value = 62520696
$ws.Range("M133").Value = 1890.9998
$ws.Range("N133").Value = -62530816
$ws.Range("H135").Value = 3476.2666
$ws.Range("I135").Value = 1168.125
$ws.Range("J135").Value = 6114.143
$ws.Range("K135").Value = 10513.125
$ws.Range("L135").Value = 55027.287
$ws.Range("M135").Value = -7978.125
$ws.Range("N135").Value = -60097.287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 459839.1
$ws.Range("I80").Value = 1002386.2
$ws.Range("J80").Value = 7716.5
$ws.Range("K80").Value = 1002386.2
$ws.Range("L80").Value = 7716.5
$ws.Range("M80").Value = -1001388.2
$ws.Range("N80").Value = -9712.5
$ws.Range("H83").Value = 459839.1
$ws.Range("I83").Value = 1002386.2
$ws.Range("J83").Value = 7716.5
$ws.Range("K83").Value = 5011931
$ws.Range("L83").Value = 38582.5
$ws.Range("M83").Value = -5006939
$ws.Range("N83").Value = -48566.5
$ws.Range("H97").Value = 4000
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 4000
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 4000
$ws.Range("N97").Value = -4992
$ws.Range("M97").ClearContents()
$ws.Range("H126").Value = 4311.3687
$ws.Range("I126").Value = 2554.111
$ws.Range("J126").Value = 5892.9
$ws.Range("K126").Value = 7662.333
$ws.Range("L126").Value = 17678.7
$ws.Range("M126").Value = -5192.333
$ws.Range("N126").Value = -22618.7
$ws.Range("H132").Value = 2950.5676
$ws.Range("I132").Value = 2364.2903
$ws.Range("J132").Value = 5979.6665
$ws.Range("K132").Value = 7092.8709
$ws.Range("L132").Value = 17938.9995
$ws.Range("M132").Value = -4562.8709
$ws.Range("N132").Value = -22998.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 11198.7
$ws.Range("I82").Value = 747.3333
$ws.Range("J82").Value = 15677.857
$ws.Range("K82").Value = 747.3333
$ws.Range("L82").Value = 15677.857
$ws.Range("M82").Value = -386.3333
$ws.Range("N82").Value = -16399.857
$ws.Range("H85").Value = 11198.7
$ws.Range("I85").Value = 747.3333
$ws.Range("J85").Value = 15677.857
$ws.Range("K85").Value = 747.3333
$ws.Range("L85").Value = 15677.857
$ws.Range("M85").Value = 500.6667
$ws.Range("N85").Value = -18173.857
$ws.Range("H93").Value = 1407.3611
$ws.Range("I93").Value = 1461.6207
$ws.Range("J93").Value = 1182.5714
$ws.Range("K93").Value = 1461.6207
$ws.Range("L93").Value = 1182.5714
$ws.Range("M93").Value = -213.6206999999999
$ws.Range("N93").Value = -3678.5714
$ws.Range("H100").Value = 66679004
$ws.Range("I100").Value = 200000000
$ws.Range("J100").Value = 18504
$ws.Range("K100").Value = 200000000
$ws.Range("L100").Value = 18504
$ws.Range("M100").Value = -199999459
$ws.Range("N100").Value = -19586
$ws.Range("H136").Value = 3367.1587
$ws.Range("I136").Value = 2656.389
$ws.Range("J136").Value = 4314.852
$ws.Range("K136").Value = 7969.167
$ws.Range("L136").Value = 12944.556
$ws.Range("M136").Value = -5419.167
$ws.Range("N136").Value = -18044.556
$ws.Range("H137").Value = 61665.4
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 61665.4
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 61665.4
$ws.Range("N137").Value = -71865.39999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2351.0264
$ws.Range("I132").Value = 1827.4706
$ws.Range("J132").Value = 6801.25
$ws.Range("K132").Value = 5482.4118
$ws.Range("L132").Value = 20403.75
$ws.Range("M132").Value = -2952.4118
$ws.Range("N132").Value = -25463.75
$ws.Range("H136").Value = 3172.3914
$ws.Range("I136").Value = 1162.8823
$ws.Range("J136").Value = 8866
$ws.Range("K136").Value = 3488.6469
$ws.Range("L136").Value = 26598
$ws.Range("M136").Value = -938.6468999999997
$ws.Range("N136").Value = -31698
